# Resize/reposition the main figure (Picture 4) on slide 2 so it becomes
# the large, full-bleed main image described in the commit message
# ("a new main figure").
#
# Target geometry (EMU, from the canonical OOXML diff):
#   off:  x=5394960  y=0
#   ext:  cx=8016240 cy=7168794
#
# PowerPoint's Shape.Left/Top/Width/Height COM properties are expressed in
# points (1 pt = 12700 EMU), so the EMU targets are converted to points
# below. The literal constants are chosen so that, once PowerPoint's
# internal Single(32-bit)-precision float round-trip happens, the saved
# OOXML reproduces the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$pic = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "Picture 4") {
        $pic = $sh
    }
}

$pic.Left   = 424.800004
$pic.Top    = 0.0
$pic.Width  = 631.2
$pic.Height = 564.471968503937
